$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The header/label row (Man'r, UNIT, QTY, U/RATE, AMOUNT, ... DATE, SOURCE,
# TYPE, PROJECT, DESCRIPTION) is removed entirely; every other row shifts
# up by one.
$ws.Rows("1:1").Delete()

# Re-assign the AMOUNT formula column as a single range write so the
# engine re-collapses it back into one shared formula group (matching the
# original authoring pattern) instead of 42 independent formula cells.
$ws.Range("F2:F43").Formula = "=D2*E2"

# Fix up the two workbook-level defined names that Excel recalculates
# when the header row disappears: the autofilter header range no longer
# exists (=> #REF!) and the print area shrinks by the one deleted row.
$fd = $wb.Names.Item("Sheet1!_FilterDatabase")
$fd.RefersTo = "=Sheet1!#REF!"

$pa = $wb.Names.Item("Sheet1!Print_Area")
$pa.RefersTo = "=Sheet1!`$A`$1:`$J`$45"

# Move the selection cursor to match where the author left off editing.
$ws.Range("E11").Select()
